# Adds a new "Volvo Assistance" automatic-crash-alert knowledge-base row
# at the top of the recently-added entries (new row 9), pushing every
# existing row down by one. Mirrors a manual Excel edit: insert a blank
# row, copy the formatting from the row above (row 8, which already has
# the "freshly added" look-and-feel) and then fill in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- insert the new row -------------------------------------------------
$ws.Rows.Item(9).Insert()

# Copy formatting (styles, borders, fill, number format, ...) from the
# row directly above so the new row matches its neighbours.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(9).RowHeight = 186

# --- fill in the new content ---------------------------------------------
$desc = @"
충돌이 발생할 경우, 차량은 Volvo Assistance 또는 긴급 콜 센터에 자동으로 보고하며, 이를 통해 긴급 지원팀이 파견될 수 있습니다.
<h3>Volvo Assistance</h3>
차량의 안전 시스템이 트리거되면(예: 안전벨트 프리텐셔너 또는 에어백 작동 레벨의 사고) 차량은 자동으로 Volvo Assistance에 전화하고 차량의 위치 등이 포함된 메시지가 전송됩니다.
<br>1. Volvo Assistance는 차량의 운전자와 통화를 하고 충돌의 정도와 도움이 필요한지 여부를 알아내려 시도합니다.
<br>2. 그런 다음 Volvo Assistance는 필요한 지원(경찰, 구급차, 차량 구조 등)을 받기 위해 연락합니다.
<br>
통화를 할 수 없는 경우에 Volvo Assistance는 적절한 조치를 통해 지원을 제공하는 관련 당국에 연락합니다.
<h3>응급 콜센터</h3><sub>대한민국은 대상 아님<sub>
차량의 안전 시스템이 트리거되면(예를 들어 안전벨트 프리텐셔너 또는 에어백이 활성화되는 수준의 사고가 발생하면) 신호가 응급 콜센터로 직접 자동 전송됩니다.
<br>
<br>1. 긴급 전화 센터는 차량의 운전자와 통화를 하고 충돌의 정도와 도움이 필요한지 여부를 알아내려 합니다.
<br>2. 긴급 전화 센터는 필요한 지원팀을 보냅니다(경찰, 구급차, 견인차 등).
"@

# Order matters for the shared-strings table: the description (D9) is
# written before the title (A9) so the new unique strings land in the
# same order as the source workbook (description first, title second).
$ws.Cells.Item(9, 4).Value = $desc
$ws.Cells.Item(9, 1).Value = "Volvo Assistance을 통한 자동 충돌 경보"
$ws.Cells.Item(9, 2).Value = 45474
$ws.Cells.Item(9, 5).Value = "Android"

# --- move the selection to D2, like the source file -----------------------
$ws.Range("D2").Select()
